# B1--and-B2-PowerPoint.pptx edit script
#
# 1) Slide 5 contains a table (graphicFrame) whose table style is switched
#    from the deck's custom "Table_0" style to the built-in
#    "Medium Style 2 - Accent 2" style (brace-GUID
#    {5578B574-13D0-4E0C-9A12-2222AF952FA3}).
#
# 2) The presentation's theme (ppt/theme/theme1.xml, used by the slide
#    master) is switched from the "Integral" / "Red Violet" palette over to
#    the stock "Office Theme" / "Office" color palette (the palette that
#    used to live alongside the notes master's theme part).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 --------------------------------------------
$s = $p.Slides.Item(5)

foreach ($shp in $s.Shapes) {
    if ($shp.HasTable) {
        $tbl = $shp.Table
        $tbl.ApplyStyle("{5578B574-13D0-4E0C-9A12-2222AF952FA3}")
    }
}

# --- 2. Theme color swap ---------------------------------------------------
# Re-point the slide master's theme color scheme at the "Office" palette
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink == ThemeColorScheme
# indexes 1-12).
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Colors(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1
$colors.Colors(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1
$colors.Colors(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2
$colors.Colors(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2
$colors.Colors(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1
$colors.Colors(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2
$colors.Colors(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3
$colors.Colors(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4
$colors.Colors(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5
$colors.Colors(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6
$colors.Colors(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink
$colors.Colors(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink
